$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the formatting used by the
# other header cells (e.g. G1 - bold/centered/bordered style).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for the "Save" column (plain, unstyled like the rest
# of the numeric data cells).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
